$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '27.564.92'
$r.ClearFormats()

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = '  +4.15%  '
$r.ClearFormats()

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '1.846.20'
$r.ClearFormats()

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '  +3.24%  '
$r.ClearFormats()

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.030'
$r.ClearFormats()

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '  +2.64%  '
$r.ClearFormats()

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '319.80'
$r.ClearFormats()

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '  +4.34%  '
$r.ClearFormats()

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '1.027'
$r.ClearFormats()

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '  +2.48%  '
$r.ClearFormats()

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '  +2.42%  '
$r.ClearFormats()

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.3745'
$r.ClearFormats()

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '  +3.66%  '
$r.ClearFormats()

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.07394'
$r.ClearFormats()

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '  +3.38%  '
$r.ClearFormats()

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.8763'
$r.ClearFormats()

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '  +2.94%  '
$r.ClearFormats()

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '21.48'
$r.ClearFormats()

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '  +4.77%  '
$r.ClearFormats()

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '1.853.23'
$r.ClearFormats()

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '  +0.88%  '
$r.ClearFormats()

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '5.488'
$r.ClearFormats()

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '  +4.23%  '
$r.ClearFormats()

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '6.689'
$r.ClearFormats()

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '  +2.99%  '
$r.ClearFormats()

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.07148'
$r.ClearFormats()

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '  +4.00%  '
$r.ClearFormats()

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '82.70'
$r.ClearFormats()

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '  +3.78%  '
$r.ClearFormats()

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '  +2.50%  '
$r.ClearFormats()

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '0.000009001'
$r.ClearFormats()

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '  +1.62%  '
$r.ClearFormats()

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = '  +2.54%  '
$r.ClearFormats()

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '15.40'
$r.ClearFormats()

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '  +2.80%  '
$r.ClearFormats()

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '27.575.98'
$r.ClearFormats()

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '  +4.13%  '
$r.ClearFormats()

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '5.264'
$r.ClearFormats()

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '  +2.49%  '
$r.ClearFormats()

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '11.20'
$r.ClearFormats()

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '  +1.81%  '
$r.ClearFormats()

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '2.064.39'
$r.ClearFormats()

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '  +0.80%  '
$r.ClearFormats()

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '157.51'
$r.ClearFormats()

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '  +3.57%  '
$r.ClearFormats()

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '1.929'
$r.ClearFormats()

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '  +5.91%  '
$r.ClearFormats()

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = '  +3.09%  '
$r.ClearFormats()

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '5.258'
$r.ClearFormats()

$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = '  +2.41%  '
$r.ClearFormats()

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '1.941'
$r.ClearFormats()

$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = '  +2.37%  '
$r.ClearFormats()

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '116.07'
$r.ClearFormats()

$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = '  +1.19%  '
$r.ClearFormats()

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '0.09083'
$r.ClearFormats()

$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = '  +2.28%  '
$r.ClearFormats()

$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = '  +3.33%  '
$r.ClearFormats()

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '1.205'
$r.ClearFormats()

$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = '  +5.59%  '
$r.ClearFormats()

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '4.498'
$r.ClearFormats()

$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = '  +3.62%  '
$r.ClearFormats()

$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = '  +4.28%  '
$r.ClearFormats()

$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = '  +2.58%  '
$r.ClearFormats()

$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = '  +2.39%  '
$r.ClearFormats()

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '0.01975'
$r.ClearFormats()

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = '  +4.14%  '
$r.ClearFormats()

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.05260'
$r.ClearFormats()

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '  +2.13%  '
$r.ClearFormats()

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.5174'
$r.ClearFormats()

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '  +4.05%  '
$r.ClearFormats()

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '2.799'
$r.ClearFormats()

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '  +7.48%  '
$r.ClearFormats()

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '  +3.25%  '
$r.ClearFormats()

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '6.676'
$r.ClearFormats()

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '  +4.27%  '
$r.ClearFormats()

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '8.513'
$r.ClearFormats()

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '  +3.67%  '
$r.ClearFormats()

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '108.90'
$r.ClearFormats()

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '  +3.29%  '
$r.ClearFormats()

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '10.55'
$r.ClearFormats()

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '  +3.05%  '
$r.ClearFormats()

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '1.713'
$r.ClearFormats()

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '  +4.65%  '
$r.ClearFormats()

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.4651'
$r.ClearFormats()

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '  +3.00%  '
$r.ClearFormats()

$r = $ws.Range("B49")
$r.NumberFormat = "@"
$r.Value = 'RenderToken'
$r.ClearFormats()

$r = $ws.Range("C49")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r.ClearFormats()

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '1.904'
$r.ClearFormats()

$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '  +8.02%  '
$r.ClearFormats()

$r = $ws.Range("B50")
$r.NumberFormat = "@"
$r.Value = 'Cronos'
$r.ClearFormats()

$r = $ws.Range("C50")
$r.NumberFormat = "@"
$r.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$r.ClearFormats()

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '0.06366'
$r.ClearFormats()

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '  +2.66%  '
$r.ClearFormats()

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '39.51'
$r.ClearFormats()

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '  +6.84%  '
$r.ClearFormats()

